$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old E1:G1 values (they were 0, 6, 4) and move them down into column D
$ws.Range("E1").ClearContents()
$ws.Range("F1").ClearContents()
$ws.Range("G1").ClearContents()

# Populate column D for rows 2-5 with the relocated values
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 6
$ws.Range("D4").Value = 4
$ws.Range("D5").Value = 1

# Update the active selection to D6
$ws.Range("D6").Select()
